$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the "amount" column values: drop the "Colombian Pesos" text, store plain numbers
$ws.Range("B1").Value = 46197879
$ws.Range("B2").Value = 44990324
$ws.Range("B3").Value = 33615976

# Apply a thousands-separator number format (numFmtId 3 = "#,##0")
$ws.Range("B1:B3").NumberFormat = "#,##0"

# Widen column B so the bigger numbers are readable
# (the engine pads ColumnWidth by 5/6 of a character when persisting to OOXML,
# so back that padding out to land on a stored width of exactly 61)
$ws.Columns.Item(2).ColumnWidth = 61 - (5/6)

# Move the active selection to B1
$ws.Range("B1").Select()
